$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Neodymium")
$ws.Range("C2").Value = [double]"6.085816639541192E-08"
$ws.Range("D2").Value = [double]"0.03496061745321188"
$ws.Range("E2").Value = [double]"0.3658596089349784"
$ws.Range("B3").Value = [double]"4.691044125953377E-15"
$ws.Range("C3").Value = [double]"0.0002268197701830409"
$ws.Range("D3").Value = [double]"0.0222307291076446"
$ws.Range("E3").Value = [double]"0.3236675139086057"
$ws.Range("B4").Value = [double]"7.321669688613381E-17"
$ws.Range("C4").Value = [double]"5.610366990890634E-05"
$ws.Range("D4").Value = [double]"0.01665716416846134"
$ws.Range("E4").Value = [double]"0.2854435934952038"
$ws.Range("C5").Value = [double]"5.113707658519867E-11"
$ws.Range("D5").Value = [double]"0.000256870885425719"
$ws.Range("E5").Value = [double]"0.02171716625600762"

$ws = $wb.Worksheets.Item("Dysprosium")
$ws.Range("C2").Value = [double]"6.085816639541599E-08"
$ws.Range("D2").Value = [double]"0.02663666091673463"
$ws.Range("E2").Value = [double]"0.3658596089350029"
$ws.Range("C3").Value = [double]"0.0002268197701830563"
$ws.Range("D3").Value = [double]"0.01693769836773036"
$ws.Range("E3").Value = [double]"0.3236675139086273"
$ws.Range("C4").Value = [double]"5.61036699089101E-05"
$ws.Range("D4").Value = [double]"0.01269117269978093"
$ws.Range("E4").Value = [double]"0.2854435934952229"
$ws.Range("C5").Value = [double]"5.113707658520278E-11"
$ws.Range("D5").Value = [double]"0.0001957111508005635"
$ws.Range("E5").Value = [double]"0.02171716625600935"

$ws = $wb.Worksheets.Item("Copper")
$ws.Range("B2").Value = [double]"3.013444709329602E-07"
$ws.Range("C2").Value = [double]"0.0001671902234504595"
$ws.Range("D2").Value = [double]"0.2022714833255239"
$ws.Range("E2").Value = [double]"0.6996950518224571"
$ws.Range("B3").Value = [double]"2.048071879604898E-06"
$ws.Range("C3").Value = [double]"0.003315520684574376"
$ws.Range("D3").Value = [double]"0.1211187136798926"
$ws.Range("E3").Value = [double]"0.4843334564770065"
$ws.Range("B4").Value = [double]"6.08104089514155E-06"
$ws.Range("C4").Value = [double]"0.0004319156678885881"
$ws.Range("D4").Value = [double]"0.07565547715308271"
$ws.Range("E4").Value = [double]"0.4657667777509991"
$ws.Range("B5").Value = [double]"1.909822742276567E-06"
$ws.Range("C5").Value = [double]"0.001007358663285219"
$ws.Range("D5").Value = [double]"0.1363938532394394"
$ws.Range("E5").Value = [double]"0.4717081075110235"

$ws = $wb.Worksheets.Item("Raw silicon")
$ws.Range("B2").Value = [double]"5.031589966768446E-08"
$ws.Range("C2").Value = [double]"5.693720412500773E-05"
$ws.Range("D2").Value = [double]"0.09118852619924153"
$ws.Range("E2").Value = [double]"1.245455460903948"
$ws.Range("B3").Value = [double]"5.367099139082082E-08"
$ws.Range("C3").Value = [double]"0.0002389556638984481"
$ws.Range("D3").Value = [double]"0.05264076056580122"
$ws.Range("E3").Value = [double]"0.5556871998096851"
$ws.Range("B4").Value = [double]"3.443534520072915E-07"
$ws.Range("C4").Value = [double]"6.181213035534014E-05"
$ws.Range("D4").Value = [double]"0.03835448965761987"
$ws.Range("E4").Value = [double]"0.5538963232362989"
$ws.Range("B5").Value = [double]"1.848515343462404E-07"
$ws.Range("C5").Value = [double]"7.794167299785755E-05"
$ws.Range("D5").Value = [double]"0.08153121125859564"
$ws.Range("E5").Value = [double]"0.883455345719077"
